# Automatic data refresh: update_automatic dades i banners [2026-02-18 04:51]
# Refreshes scraped meteo.cat observations (extraction timestamps + measured
# values) for rows 2-46 of the Dades_Meteo sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Humidity cells hold plain-text percentages ("70%"). Excel's literal-input
# parser would otherwise read "NN%" as the number NN/100, so force the Text
# number format on these specific cells before writing the new value.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "70%"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "91%"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "95%"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "96%"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "94%"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "93%"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "94%"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "84%"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "83%"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "55%"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "90%"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "75%"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "84%"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "57%"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "79%"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "91%"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "38%"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "97%"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "69%"

# Remaining updated cells: extraction timestamps and measurement text that
# already round-trip as literal strings (units/symbols prevent numeric parsing).
$ws.Range("E2").Value = "2026-02-18 04:49:04"
$ws.Range("O2").Value = "-0.9 °C"
$ws.Range("E3").Value = "2026-02-18 04:49:06"
$ws.Range("M3").Value = "-2.7 °C 4:29 TU"
$ws.Range("O3").Value = "-3.8 °C"
$ws.Range("E4").Value = "2026-02-18 04:49:09"
$ws.Range("J4").Value = "1018.7 hPa"
$ws.Range("N4").Value = "6.5 °C 4:21 TU"
$ws.Range("O4").Value = "7.6 °C"
$ws.Range("E5").Value = "2026-02-18 04:49:12"
$ws.Range("K5").Value = "-0.1 MJ/m2"
$ws.Range("M5").Value = "-0.9 °C 4:25 TU"
$ws.Range("O5").Value = "-2.4 °C"
$ws.Range("E6").Value = "2026-02-18 04:49:14"
$ws.Range("J6").Value = "1018.4 hPa"
$ws.Range("N6").Value = "6.8 °C 4:14 TU"
$ws.Range("O6").Value = "8.4 °C"
$ws.Range("E7").Value = "2026-02-18 04:49:17"
$ws.Range("J7").Value = "1018.6 hPa"
$ws.Range("N7").Value = "12.1 °C 4:26 TU"
$ws.Range("E8").Value = "2026-02-18 04:49:19"
$ws.Range("J8").Value = "1018.7 hPa"
$ws.Range("N8").Value = "8.1 °C 4:24 TU"
$ws.Range("O8").Value = "9.0 °C"
$ws.Range("E9").Value = "2026-02-18 04:49:22"
$ws.Range("L9").Value = "9.7 km/h - 298º 4:16 TU"
$ws.Range("E10").Value = "2026-02-18 04:49:24"
$ws.Range("N10").Value = "4.9 °C 4:29 TU"
$ws.Range("O10").Value = "7.5 °C"
$ws.Range("E11").Value = "2026-02-18 04:49:27"
$ws.Range("N11").Value = "0.5 °C 4:29 TU"
$ws.Range("O11").Value = "2.6 °C"
$ws.Range("E12").Value = "2026-02-18 04:49:30"
$ws.Range("O12").Value = "6.3 °C"
$ws.Range("E13").Value = "2026-02-18 04:49:32"
$ws.Range("O13").Value = "-1.9 °C"
$ws.Range("E14").Value = "2026-02-18 04:49:35"
$ws.Range("N14").Value = "8.4 °C 4:28 TU"
$ws.Range("O14").Value = "10.4 °C"
$ws.Range("E15").Value = "2026-02-18 04:49:37"
$ws.Range("O15").Value = "5.4 °C"
$ws.Range("E16").Value = "2026-02-18 04:49:40"
$ws.Range("E17").Value = "2026-02-18 04:49:43"
$ws.Range("O17").Value = "2.0 °C"
$ws.Range("E18").Value = "2026-02-18 04:49:45"
$ws.Range("J18").Value = "1018.7 hPa"
$ws.Range("N18").Value = "5.7 °C 4:26 TU"
$ws.Range("O18").Value = "8.1 °C"
$ws.Range("E19").Value = "2026-02-18 04:49:48"
$ws.Range("L19").Value = "9.4 km/h - 245º 4:19 TU"
$ws.Range("E20").Value = "2026-02-18 04:49:51"
$ws.Range("L20").Value = "30.2 km/h - 284º 4:12 TU"
$ws.Range("O20").Value = "-1.3 °C"
$ws.Range("E21").Value = "2026-02-18 04:49:53"
$ws.Range("N21").Value = "0.4 °C 4:10 TU"
$ws.Range("O21").Value = "2.2 °C"
$ws.Range("E22").Value = "2026-02-18 04:49:56"
$ws.Range("E23").Value = "2026-02-18 04:49:58"
$ws.Range("O23").Value = "0.2 °C"
$ws.Range("E24").Value = "2026-02-18 04:50:01"
$ws.Range("J24").Value = "1019.1 hPa"
$ws.Range("N24").Value = "3.2 °C 4:17 TU"
$ws.Range("O24").Value = "5.8 °C"
$ws.Range("E25").Value = "2026-02-18 04:50:04"
$ws.Range("O25").Value = "-0.6 °C"
$ws.Range("E26").Value = "2026-02-18 04:50:06"
$ws.Range("E27").Value = "2026-02-18 04:50:09"
$ws.Range("L27").Value = "27.4 km/h - 233º 4:29 TU"
$ws.Range("E28").Value = "2026-02-18 04:50:11"
$ws.Range("J28").Value = "1019.1 hPa"
$ws.Range("O28").Value = "5.3 °C"
$ws.Range("E29").Value = "2026-02-18 04:50:13"
$ws.Range("N29").Value = "7.9 °C 4:29 TU"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-18 04:50:16"
$ws.Range("J30").Value = "1018.7 hPa"
$ws.Range("E31").Value = "2026-02-18 04:50:19"
$ws.Range("J31").Value = "1017.4 hPa"
$ws.Range("N31").Value = "10.0 °C 4:28 TU"
$ws.Range("O31").Value = "10.7 °C"
$ws.Range("E32").Value = "2026-02-18 04:50:21"
$ws.Range("L32").Value = "6.8 km/h - 318º 4:25 TU"
$ws.Range("E33").Value = "2026-02-18 04:50:24"
$ws.Range("J33").Value = "1021.4 hPa"
$ws.Range("N33").Value = "-1.0 °C 4:28 TU"
$ws.Range("O33").Value = "0.2 °C"
$ws.Range("E34").Value = "2026-02-18 04:50:26"
$ws.Range("O34").Value = "-0.7 °C"
$ws.Range("E35").Value = "2026-02-18 04:50:29"
$ws.Range("N35").Value = "6.4 °C 4:19 TU"
$ws.Range("O35").Value = "7.7 °C"
$ws.Range("E36").Value = "2026-02-18 04:50:32"
$ws.Range("E37").Value = "2026-02-18 04:50:35"
$ws.Range("J37").Value = "1021.4 hPa"
$ws.Range("E38").Value = "2026-02-18 04:50:37"
$ws.Range("N38").Value = "7.4 °C 4:29 TU"
$ws.Range("O38").Value = "9.5 °C"
$ws.Range("E39").Value = "2026-02-18 04:50:39"
$ws.Range("L39").Value = "38.2 km/h - 247º 4:15 TU"
$ws.Range("M39").Value = "1.7 °C 4:06 TU"
$ws.Range("O39").Value = "-0.1 °C"
$ws.Range("E40").Value = "2026-02-18 04:50:42"
$ws.Range("N40").Value = "-0.4 °C 4:29 TU"
$ws.Range("O40").Value = "1.0 °C"
$ws.Range("E41").Value = "2026-02-18 04:50:44"
$ws.Range("J41").Value = "1018.3 hPa"
$ws.Range("N41").Value = "7.2 °C 4:29 TU"
$ws.Range("E42").Value = "2026-02-18 04:50:47"
$ws.Range("N42").Value = "6.4 °C 4:07 TU"
$ws.Range("O42").Value = "8.7 °C"
$ws.Range("E43").Value = "2026-02-18 04:50:49"
$ws.Range("N43").Value = "6.2 °C 4:19 TU"
$ws.Range("O43").Value = "7.3 °C"
$ws.Range("E44").Value = "2026-02-18 04:50:52"
$ws.Range("O44").Value = "-3.7 °C"
$ws.Range("E45").Value = "2026-02-18 04:50:54"
$ws.Range("J45").Value = "1021.0 hPa"
$ws.Range("L45").Value = "10.4 km/h - 169º 4:25 TU"
$ws.Range("E46").Value = "2026-02-18 04:50:57"
$ws.Range("J46").Value = "1019.2 hPa"
$ws.Range("N46").Value = "4.8 °C 4:29 TU"
$ws.Range("O46").Value = "6.7 °C"
